# Auto-generated edit script: updates the cryptos price/volume table
# to match the refreshed GitHub Actions data pull, and swaps the
# Aave / BabyDogeCoin rows (45 <-> 46) to their new ranking order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.826.80"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").Value = "1.599.26"
$ws.Range("E3").Value = "  -2.11%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'208.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.50%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "'0.477"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.37%  "
$ws.Range("D8").Value = "'0.246"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.59%  "
$ws.Range("D9").Value = "'0.0610"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.07%  "
$ws.Range("D10").Value = "'17.86"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.42%  "
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("D12").Value = "1.819.20"
$ws.Range("E12").Value = "  -2.23%  "
$ws.Range("D13").Value = "1.603.74"
$ws.Range("E13").Value = "  -1.98%  "
$ws.Range("E14").Value = "  -3.53%  "
$ws.Range("D15").Value = "'0.508"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "25.819.62"
$ws.Range("E16").Value = "  -0.63%  "
$ws.Range("E17").Value = "  -2.06%  "
$ws.Range("D18").Value = "0.0₃0716"
$ws.Range("E18").Value = "  -4.07%  "
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").Value = "'189.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("E21").Value = "  -1.58%  "
$ws.Range("D22").Value = "'9.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.78%  "
$ws.Range("E23").Value = "  -3.29%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("E25").Value = "  -3.04%  "
$ws.Range("D26").Value = "'141.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.33%  "
$ws.Range("D27").Value = "'1.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.16%  "
$ws.Range("D28").Value = "'6.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.18%  "
$ws.Range("D29").Value = "'14.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.19%  "
$ws.Range("E30").Value = "  -2.83%  "
$ws.Range("D31").Value = "'0.0461"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.82%  "
$ws.Range("D32").Value = "'3.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.42%  "
$ws.Range("E33").Value = "  -5.03%  "
$ws.Range("E34").Value = "  -1.09%  "
$ws.Range("E35").Value = "  -2.04%  "
$ws.Range("D36").Value = "1.098.31"
$ws.Range("E36").Value = "  -3.20%  "
$ws.Range("E37").Value = "  -3.19%  "
$ws.Range("E38").Value = "  -7.95%  "
$ws.Range("D39").Value = "'0.0150"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.82%  "
$ws.Range("D40").Value = "'0.493"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.99%  "
$ws.Range("D41").Value = "'95.66"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.94%  "
$ws.Range("D42").Value = "1.732.88"
$ws.Range("E42").Value = "  -2.13%  "
$ws.Range("D43").Value = "'5.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.74%  "
$ws.Range("D44").Value = "'0.741"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.99%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'53.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.76%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₇0996"
$ws.Range("E46").Value = "  -12.71%  "
$ws.Range("E47").Value = "  -3.20%  "
$ws.Range("E48").Value = "  -3.57%  "
$ws.Range("E49").Value = "  -0.83%  "
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("D51").Value = "'7.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.60%  "
